$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Controllers")

# Volt-Var debugging: swap the TS1/TS2 pyPlotList time-series entries to the
# new device (261584/261585 instead of 263163/263265) and switch the TS2
# index filter from the single-index form to the even-sample form.
$ws.Range("B3").Value = "PVSystem.oh_261584_2_4"
$ws.Range("B4").Value = "dev_261585_2_4"
$ws.Range("F4").Value = "['Even']"

# The TS1/TS2 rows were manually resized while reviewing the settings.
$ws.Rows.Item(3).RowHeight = 18.75
$ws.Rows.Item(4).RowHeight = 20.25

# Leave the cursor parked on the cell that was being edited (F4).
$ws.Activate() | Out-Null
$ws.Range("F4").Select() | Out-Null
